# Refresh the cryptos price/volume(1h) columns with the latest scraped
# values (GitHub Actions scheduled update, Thu Nov 7 2024 20:43:48 UTC).
# Price (D) and Volume(1h) (E) cells are stored as literal text in this
# sheet, so NumberFormat is forced to "@" (text) before each Price write
# to stop Excel from silently re-interpreting strings such as "1.00" or
# "0.0000190" as numbers and dropping their trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '76.470.54'
$ws.Range('E2').Value = '  +1.17%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.886.21'
$ws.Range('E3').Value = '  +7.96%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '197.02'
$ws.Range('E5').Value = '  +5.33%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '600.21'
$ws.Range('E6').Value = '  +2.31%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  +3.21%  '
$ws.Range('E9').Value = '  -0.45%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '2.883.93'
$ws.Range('E10').Value = '  +7.88%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.399'
$ws.Range('E11').Value = '  +11.37%  '
$ws.Range('E12').Value = '  -1.73%  '
$ws.Range('E13').Value = '  +4.26%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.391.02'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '76.383.94'
$ws.Range('E15').Value = '  +1.28%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '27.54'
$ws.Range('E16').Value = '  +3.90%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000190'
$ws.Range('E17').Value = '  +0.93%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.889.69'
$ws.Range('E18').Value = '  +7.85%  '
$ws.Range('E19').Value = '  -3.07%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.61'
$ws.Range('E20').Value = '  +5.56%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '383.44'
$ws.Range('E21').Value = '  +3.00%  '
$ws.Range('E22').Value = '  +1.55%  '
$ws.Range('E23').Value = '  +1.57%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '71.91'
$ws.Range('E24').Value = '  +2.81%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.00'
$ws.Range('E25').Value = '  +0.05%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.032.67'
$ws.Range('E26').Value = '  +7.29%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '4.24'
$ws.Range('E27').Value = '  +1.49%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.82'
$ws.Range('E28').Value = '  +4.90%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0000106'
$ws.Range('E29').Value = '  +12.08%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.997'
$ws.Range('E30').Value = '  -0.40%  '
$ws.Range('E31').Value = '  +1.07%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '514.31'
$ws.Range('E32').Value = '  -1.10%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.79'
$ws.Range('E33').Value = '  +0.68%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.82'
$ws.Range('E34').Value = '  +3.83%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '168.09'
$ws.Range('E36').Value = '  +2.86%  '
$ws.Range('E37').Value = '  +4.94%  '
$ws.Range('E38').Value = '  -1.74%  '
$ws.Range('E39').Value = '  +0.83%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '182.68'
$ws.Range('E40').Value = '  +8.38%  '
$ws.Range('E41').Value = '  -0.07%  '
$ws.Range('E42').Value = '  +5.21%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.10'
$ws.Range('E43').Value = '  +1.89%  '
$ws.Range('E44').Value = '  -0.43%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0931'
$ws.Range('E45').Value = '  +10.30%  '
$ws.Range('E46').Value = '  +4.00%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '40.21'
$ws.Range('E47').Value = '  +2.23%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.36'
$ws.Range('E48').Value = '  +0.61%  '
$ws.Range('E49').Value = '  +18.12%  '
$ws.Range('E50').Value = '  +8.66%  '
$ws.Range('E51').Value = '  +3.49%  '
